# Update MSME Country Indicators - Ukraine Summary data values.
# These cells hold numeric-looking text (shared strings), so we assign
# apostrophe-prefixed values to force Excel to keep them as text rather
# than auto-converting to numbers, then reset the cell style back to
# Normal so no quote-prefix formatting lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, $text) {
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# Enterprises density (per 1000 people) - row 11
Set-TextValue "B11" "6.28"
Set-TextValue "C11" "1.71"
Set-TextValue "D11" "7.99"

# Employment (% of total) - row 12
Set-TextValue "B12" "10.26"
Set-TextValue "C12" "57.39"
Set-TextValue "D12" "67.65"

# Enterprises (% of total) - row 14 (only C and D change)
Set-TextValue "C14" "21.31"
Set-TextValue "D14" "99.81"
